# Change the year in the astromap link: GaNight/2018 -> GaNight/2022.
# The old text is a HYPERLINK field (begin/instrText/separate/result/end)
# wrapped in "(" and ")." plain-text runs. The new text is plain (unlinked)
# text: "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$d = $word.ActiveDocument

# Locate the paragraph holding the old link by searching for its
# (still-)unique year-stamped path fragment.
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("GaNight/2018", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)

if ($found) {
    $hit = $find.Parent
    $hit.Expand(4)            # wdParagraph -> whole paragraph
    $hit.MoveEnd(1, -1)       # drop the trailing paragraph mark

    # Wipe the paragraph's content (the "(" run, the hyperlink field, and
    # the ")." run) and replace it with a single plain-text run.
    $hit.Delete()
    $hit.InsertAfter("(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).")
}
